# LOM3106.docx restructuring.
#
# The body text of several sections gets rotated to a different heading
# (paragraph styles / bold labels stay exactly where they are -- only the
# "value" text moves). Concretely (old location -> new location):
#
#   Objetivos (PT)          -> Programa resumido (PT)
#   Objetivos (EN, italic)  -> Programa resumido (EN, italic)
#   Docente(s) list         -> split: first name -> Avaliacao/Norma de recuperacao
#                                      second name -> Bibliografia
#   Programa resumido (PT)  -> Avaliacao/Metodo
#   Programa resumido (EN)  -> Objetivos (EN, italic)
#   Programa (PT list)      -> Avaliacao/Criterio
#   Avaliacao/Metodo        -> Avaliacao/Norma de recuperacao... (chained, see below)
#   Avaliacao/Criterio      -> Bibliografia
#   Avaliacao/Norma         -> Docente(s) (first bullet)
#   Bibliografia            -> Docente(s) (second bullet)
#   (Objetivos(PT) + Programa(PT list)) together become the new Docente(s)-styled
#     paragraph's content (that paragraph keeps the ListBullet style but now
#     holds the old Objetivos/Programa text instead of the names).
#
# To apply this safely with Word's Find/Replace (which only sees flat text,
# not which logical "block" a string came from) we do it in two passes:
#   1) stamp every source block with a unique, never-colliding placeholder
#   2) resolve every placeholder into its final text
# That way a later rule can never accidentally re-match text that an earlier
# rule just inserted.

$d = $word.ActiveDocument
$vt = [char]11   # manual line break char -- round-trips through <w:br/>

function Replace-Exact($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        throw "Find.Execute could not find: $find"
    }
}

# ---------------------------------------------------------------------------
# Source texts (exactly as they appear in the original document)
# ---------------------------------------------------------------------------

$objetivosPT = "Possibilitar ao estudante de Engenharia de Materiais o acesso a ferramentas computacionais modernas, de modo a que consiga descrever e quantificar conceitos vistos em outras disciplinas, como Ciência dos Materiais, Diagramas de Fases, Cinética de Transformação em Materiais, Termodinâmica, Propriedades Elétricas, Magnéticas, Térmicas e Ópticas, etc. Ao final do curso, o aluno será capaz de aplicar e entender resultados de simulações computacionais realistas aplicadas a diversas classes de materiais."

$objetivosEN = "Provide to Materials Engineering students access to modern computational tools, so that they can describe and quantify concepts seen in other disciplines, such as Materials Science, Phase Diagrams, Transformation Kinetics in Materials, Thermodynamics, Electrical, Magnetic, Thermal and Optical Properties, etc. At the end of the course, the student will be able to apply and understand the results of realistic computer simulations applied to different classes of materials."

$docenteLine1 = "3480026 - João Paulo Pascon"
$docenteLine2 = "1176388 - Luiz Tadeu Fernandes Eleno"
$docentesFull = $docenteLine1 + $vt + $docenteLine2

$progResumidoPT = "Tratamento de imagens em materialografia; Ajuste de equações empíricas ; Potenciais interatômicos e dinâmica molecular clássica; Descrição da Cinética de nucleação e crescimento; Método dos Elementos Finitos; Métodos de Monte Carlo; Crescimento de grão; Cálculo de Diagramas de fases."

$progResumidoEN = "Image processing in materialography; Adjusting empirical equations; Interatomic potentials and classical molecular dynamics; Description of nucleation and growth kinetics; Finite Element Method; Monte Carlo methods; Grain growth; Calculation of phase diagrams."

$programaPT = "- Tratamento de imagens: resolução, definição, contraste, saturação; uso de técnicas automatizadas de determinação de tamanho e distribuição de partículas." + $vt + "- Proposição e ajuste de equações empíricas a resultados de medidas experimentais: as diversas propostas de relações para a deformação plástica e encruamento." + $vt + "- Potenciais interatômicos e o método de dinâmica molecular clássica; simulação de solidificação de um metal puro." + $vt + "- Cinética de nucleação e crescimento: a equação de Johnson-Mehl-Avrami-Kolmogorov (JMAK) e sua aplicação computacional." + $vt + "- Elementos finitos: estudo do estado de tensão de materiais sob carregamentos mecânicos; simulação de transferência de calor em tratamentos térmicos." + $vt + "- Método de Monte Carlo aplicado à transição ferro-paramagnética e à cinética de crescimento de grão" + $vt + "- Cálculo de diagramas de fases: curvas de energia livre, o método CALPHAD; Thermo-Calc e Dictra."

$metodoVal = "Aulas expositivas e em laboratório computacional, trabalhos e exercícios comentados. Trabalho baseado em Projeto"
$criterioVal = "Média aritmética de trabalhos propostos ao longo do curso (60%) e do Trabalho final em grupo (40%)."
$normaVal = "Não haverá exame de recuperação."

$bibliografia = "- Richard LESAR, Computational Materials Science – Fundamentals to Applications. MRS, 2013." + $vt + "- Rob Phillips, Crystals, Defects and Microstructures – Modelling across scales. Cambridge, 2001." + $vt + "- Artigos publicados em revistas como Computational Materials Science, Calphad, Journal of Alloys and Compounds, etc."

# ---------------------------------------------------------------------------
# Phase 1: stamp each source block with a unique placeholder token.
# The docente list is stamped as a whole (it gets consumed entirely by the
# paragraph-9 rebuild below), its two names are produced later verbatim from
# the literal strings above.
# ---------------------------------------------------------------------------

Replace-Exact $objetivosPT     "@@OBJETIVOS_PT@@"
Replace-Exact $objetivosEN     "@@OBJETIVOS_EN@@"
Replace-Exact $docentesFull    "@@DOCENTES_OLD@@"
Replace-Exact $progResumidoPT  "@@PRESUMIDO_PT@@"
Replace-Exact $progResumidoEN  "@@PRESUMIDO_EN@@"
Replace-Exact $programaPT      "@@PROGRAMA_PT@@"
Replace-Exact $metodoVal       "@@METODO_VAL@@"
Replace-Exact $criterioVal     "@@CRITERIO_VAL@@"
Replace-Exact $normaVal        "@@NORMA_VAL@@"
Replace-Exact $bibliografia    "@@BIBLIOGRAFIA@@"

# ---------------------------------------------------------------------------
# Phase 2: resolve placeholders into final text.
# ---------------------------------------------------------------------------

# Objetivos section body -> old "Programa resumido" text
Replace-Exact "@@OBJETIVOS_PT@@" $progResumidoPT
Replace-Exact "@@OBJETIVOS_EN@@" $progResumidoEN

# Programa resumido section body -> old Avaliacao/Metodo (PT) and old Objetivos (EN)
Replace-Exact "@@PRESUMIDO_PT@@" $metodoVal
Replace-Exact "@@PRESUMIDO_EN@@" $objetivosEN

# Programa section body -> old Avaliacao/Criterio value
Replace-Exact "@@PROGRAMA_PT@@" $criterioVal

# Avaliacao values -> rotate
Replace-Exact "@@METODO_VAL@@"   $normaVal
Replace-Exact "@@CRITERIO_VAL@@" $bibliografia
Replace-Exact "@@NORMA_VAL@@"    $docenteLine1

# Bibliografia body -> old second docente line
Replace-Exact "@@BIBLIOGRAFIA@@" $docenteLine2

# Docente(s) paragraph (old "@@DOCENTES_OLD@@" placeholder) -> becomes the old
# Programa (PT list) body; then the old Objetivos (PT) text + a line break is
# inserted in front of it, recreating the two-run paragraph from the target.
Replace-Exact "@@DOCENTES_OLD@@" $programaPT

$docenteParaIndex = 9
$p9 = $d.Paragraphs.Item($docenteParaIndex)
if ($p9.Range.Text -notmatch [regex]::Escape($programaPT.Substring(0, 30))) {
    throw "Paragraph 9 does not contain the expected Programa text; aborting."
}
$p9.Range.InsertBefore($objetivosPT + $vt)

Write-Output "LOM3106 restructuring complete."
